# Rename the two Pearson logo pictures (in the first-page footer and the
# default footer) from "image2.png" to "image1.png", and rename the BTec
# logo picture (in the first-page header) from "image1.jpg" to "image2.jpg".
#
# These logos are inline pictures that live in headers/footers rather than
# the main document body, so we walk Sections -> Headers/Footers ->
# Range.InlineShapes and match each picture by its (unchanged) alt text
# description so the right new name lands on the right picture regardless
# of collection ordering.

$d = $word.ActiveDocument

$renameMap = @{
    "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" = "image1.png"
    "BTec_Logo-Orange" = "image2.jpg"
}

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $sec = $d.Sections.Item($s)

    for ($h = 1; $h -le $sec.Headers.Count; $h++) {
        $hdr = $sec.Headers.Item($h)
        if ($hdr.Exists) {
            for ($i = 1; $i -le $hdr.Range.InlineShapes.Count; $i++) {
                $shp = $hdr.Range.InlineShapes.Item($i)
                $newName = $renameMap[$shp.AlternativeText]
                if ($newName) {
                    $shp.Name = $newName
                }
            }
        }
    }

    for ($f = 1; $f -le $sec.Footers.Count; $f++) {
        $ftr = $sec.Footers.Item($f)
        if ($ftr.Exists) {
            for ($i = 1; $i -le $ftr.Range.InlineShapes.Count; $i++) {
                $shp = $ftr.Range.InlineShapes.Item($i)
                $newName = $renameMap[$shp.AlternativeText]
                if ($newName) {
                    $shp.Name = $newName
                }
            }
        }
    }
}
